$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B45 should become a true numeric value (it was a text "2" before)
$ws.Cells.Item(45, 2).Value = 2

# Append new row 46 with the new annotation data
$ws.Cells.Item(46, 1).Value = "Sunsi Wu"
$ws.Cells.Item(46, 2).Value = "'3"
$ws.Cells.Item(46, 2).Style = "Normal"
$ws.Cells.Item(46, 3).Value = "无"
$ws.Cells.Item(46, 4).Value = "SMY"
$ws.Cells.Item(46, 5).Value = "RES"
$ws.Cells.Item(46, 6).Value = "18e2478f-5f8b-460a-bbaf-4b86b95999fd"
$ws.Cells.Item(46, 7).Value = "B1IDRdeCW_annotated.xlsx"
$ws.Cells.Item(46, 8).Value = "This paper presents three observations to understand binary network in Courbariaux, Hubara et al. (2016)."
